$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp2 = $s.Shapes.Item(2)
$shp2.Delete()
$shp2b = $s.Shapes.Item(2)
$tr = $shp2b.TextFrame.TextRange
$tr.Text = "seed"
$tr.Text = "(Fuck This Protocol)`rby William Orr"
$para1 = $tr.Paragraphs(1,1)
$para1.Font.Name = "Arial"
